# Deploy the implementation guide.
#
# 1. Metadata!B6  "active" -> "draft"            (Status value)
# 2. Metadata!B8  old timestamp -> new timestamp  (Date value)
# 3. Re-assert wrap-text alignment on the two bordered cell styles used
#    throughout both sheets (header row + body rows) so the alignment
#    formatting ("applyAlignment"/wrapText) is (re)applied.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Metadata")
$ws2 = $wb.Worksheets.Item("Include from Ferlab.bio CodeS")

# --- Value edits -----------------------------------------------------
$ws1.Range("B6").Value = "draft"
$ws1.Range("B8").Value = "2023-08-01T16:12:28+00:00"

# --- Formatting edits --------------------------------------------------
# Header row (bold / filled / bordered style) on both sheets.
$ws1.Range("A1:B1").WrapText = $true
$ws2.Range("A1").WrapText = $true

# Body rows (plain bordered style) on both sheets - only touch cells
# that already exist/carry the style so no new blank cells are added.
$ws1.Range("A2:B14").WrapText = $true
$ws2.Range("A2:A4").WrapText = $true
$ws2.Range("B3:B4").WrapText = $true
